# Update column G ("K") values on Sheet1 per regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 1
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 4
    17 = 2
    18 = 2
    19 = 2
    21 = 1
    22 = 2
    24 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
